# Fixed more column size adjusting issues
#
# Widen the 6 data columns on the "Pastry Sales Data" sheet to their new,
# slightly larger auto-fit widths.
#
# Excel's ColumnWidth setter (and this host's xlsx writer) snaps the stored
# OOXML <col width="..."> to the nearest 1/6-character pixel-grid value via
#   stored = (Round(ColumnWidth * 6) + 5) / 6
# so we pre-compensate each requested width by subtracting the fixed 5/6
# offset before assigning it, landing as close as the grid allows on the
# intended final widths below:
#   A -> 15.025425
#   B -> 14.335425
#   C -> 24.595425
#   D -> 16.155425
#   E -> 15.685425
#   F -> 18.055425

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pastry Sales Data")

$targetWidths = @(15.025425, 14.335425, 24.595425, 16.155425, 15.685425, 18.055425)
$offset = 5.0 / 6.0

for ($i = 0; $i -lt $targetWidths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $targetWidths[$i] - $offset
}
